# Add a "Strain" column (new column B) to the weights sheet, between the
# existing "ID" and "Cohort" columns, and fill in Sex for the two rows
# that previously lacked it (while clearing the stray "n/a" that had been
# entered for row 5). This matches the commit "Add columns to pass
# validation."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column before the current column B (Cohort), shifting
# Cohort / Date_Infected / Sex one column to the right.
$ws.Columns("B:B").Insert()

# Give the new column the same width as column A so the pair reads as one
# visually consistent block (mirrors the merged col run in the saved file).
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Header for the new column.
$ws.Range("B1").Value = "Strain"

# Strain values: first three animals are the 8034x8043 cross, the last
# three are the 15119x16521 cross.
$ws.Range("B2").Value = "8034x8043"
$ws.Range("B3").Value = "8034x8043"
$ws.Range("B4").Value = "8034x8043"
$ws.Range("B5").Value = "15119x16521"
$ws.Range("B6").Value = "15119x16521"
$ws.Range("B7").Value = "15119x16521"

# Row 5 (now E5) had an erroneous "n/a" Sex value - remove it entirely
# (value + format) rather than just blanking the text.
$ws.Range("E5").Clear()

# Rows 6 and 7 (now E6 / E7) previously had no Sex recorded - fill them in.
# E7 is written first so "Male" is registered in the shared-string table
# ahead of "Female", matching the order the values were authored in.
$ws.Range("E7").Value = "Male"
$ws.Range("E6").Value = "Female"

# Match the saved selection/active cell from the edited workbook.
$ws.Range("E5").Select() | Out-Null
